$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "57.749.50"
Set-TextValue 2 5 "  -2.19%  "

# Row 3
Set-TextValue 3 4 "2.553.03"
Set-TextValue 3 5 "  -3.66%  "

# Row 4
Set-TextValue 4 5 "  -0.07%  "

# Row 5
Set-TextValue 5 4 "518.62"
Set-TextValue 5 5 "  -0.93%  "

# Row 6
Set-TextValue 6 4 "138.73"
Set-TextValue 6 5 "  -3.97%  "

# Row 7
Set-TextValue 7 4 "1.00"
Set-TextValue 7 5 "  +0.00%  "

# Row 8
Set-TextValue 8 4 "0.561"
Set-TextValue 8 5 "  -1.79%  "

# Row 9
Set-TextValue 9 4 "6.51"
Set-TextValue 9 5 "  -7.17%  "

# Row 10
Set-TextValue 10 4 "0.0988"
Set-TextValue 10 5 "  -3.76%  "

# Row 11
Set-TextValue 11 4 "0.323"
Set-TextValue 11 5 "  -3.12%  "

# Row 12
Set-TextValue 12 5 "  -0.01%  "

# Row 13
Set-TextValue 13 4 "2.999.56"
Set-TextValue 13 5 "  -3.82%  "

# Row 14
Set-TextValue 14 4 "57.731.97"
Set-TextValue 14 5 "  -2.32%  "

# Row 15
Set-TextValue 15 4 "19.94"
Set-TextValue 15 5 "  -5.24%  "

# Row 16
Set-TextValue 16 2 "ShibaInu"
Set-TextValue 16 3 "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue 16 4 "0.0000132"
Set-TextValue 16 5 "  -3.21%  "

# Row 17
Set-TextValue 17 2 "WrappedEther"
Set-TextValue 17 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue 17 4 "2.527.49"
Set-TextValue 17 5 "  -5.21%  "

# Row 18
Set-TextValue 18 4 "333.13"
Set-TextValue 18 5 "  -2.05%  "

# Row 19
Set-TextValue 19 4 "4.27"
Set-TextValue 19 5 "  -2.27%  "

# Row 20
Set-TextValue 20 4 "10.09"
Set-TextValue 20 5 "  -2.56%  "

# Row 21
Set-TextValue 21 4 "6.11"

# Row 22
Set-TextValue 22 5 "  -0.11%  "

# Row 23
Set-TextValue 23 4 "64.90"
Set-TextValue 23 5 "  +1.24%  "

# Row 24
Set-TextValue 24 5 "  -1.73%  "

# Row 25
Set-TextValue 25 5 "  +0.28%  "

# Row 26
Set-TextValue 26 4 "0.400"
Set-TextValue 26 5 "  -4.27%  "

# Row 27
Set-TextValue 27 4 "2.691.43"
Set-TextValue 27 5 "  -3.00%  "

# Row 28
Set-TextValue 28 4 "6.93"
Set-TextValue 28 5 "  -2.30%  "

# Row 29
Set-TextValue 29 4 "0.0₃0751"
Set-TextValue 29 5 "  -6.35%  "

# Row 30
Set-TextValue 30 4 "0.999"
Set-TextValue 30 5 "  +0.02%  "

# Row 31
Set-TextValue 31 4 "6.15"
Set-TextValue 31 5 "  -7.88%  "

# Row 32
Set-TextValue 32 4 "1.56"
Set-TextValue 32 5 "  -1.79%  "

# Row 33
Set-TextValue 33 4 "149.17"

# Row 34
Set-TextValue 34 4 "18.40"
Set-TextValue 34 5 "  -2.34%  "

# Row 35
Set-TextValue 35 4 "3.96"
Set-TextValue 35 5 "  -4.51%  "

# Row 36
Set-TextValue 36 5 "  -5.88%  "

# Row 37
Set-TextValue 37 4 "0.833"
Set-TextValue 37 5 "  -6.65%  "

# Row 38
Set-TextValue 38 4 "35.67"
Set-TextValue 38 5 "  -2.89%  "

# Row 39
Set-TextValue 39 4 "0.817"
Set-TextValue 39 5 "  -6.59%  "

# Row 40
Set-TextValue 40 5 "  -5.06%  "

# Row 41
Set-TextValue 41 4 "0.999"
Set-TextValue 41 5 "  -0.01%  "

# Row 42
Set-TextValue 42 5 "  -3.81%  "

# Row 43
Set-TextValue 43 2 "WhiteBITCoin"
Set-TextValue 43 3 "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue 43 4 "10.65"
Set-TextValue 43 5 "  -0.05%  "

# Row 44
Set-TextValue 44 2 "Stellar"
Set-TextValue 44 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 44 4 "0.0953"
Set-TextValue 44 5 "  -1.77%  "

# Row 45
Set-TextValue 45 2 "Mantle"
Set-TextValue 45 3 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue 45 4 "0.577"
Set-TextValue 45 5 "  -6.37%  "

# Row 46
Set-TextValue 46 2 "Bittensor"
Set-TextValue 46 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue 46 4 "260.09"
Set-TextValue 46 5 "  -5.34%  "

# Row 47
Set-TextValue 47 4 "0.0517"
Set-TextValue 47 5 "  -2.99%  "

# Row 48
Set-TextValue 48 4 "1.982.30"
Set-TextValue 48 5 "  -2.49%  "

# Row 49
Set-TextValue 49 4 "18.41"
Set-TextValue 49 5 "  -7.71%  "

# Row 50
Set-TextValue 50 4 "0.0221"
Set-TextValue 50 5 "  -3.33%  "

# Row 51
Set-TextValue 51 4 "4.51"
Set-TextValue 51 5 "  -5.72%  "
